$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.351.60"
$ws.Range("E2").Value = "  -2.61%  "
$ws.Range("D3").Value = "1.942.36"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.29"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7221"
$ws.Range("E6").Value = "  -8.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3352"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.82"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07399"
$ws.Range("E10").Value = "  +5.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8203"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08146"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.941.39"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.506"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.41"
$ws.Range("E15").Value = "  -5.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.93"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "30.374.10"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008364"
$ws.Range("E18").Value = "  +5.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.14"
$ws.Range("E19").Value = "  -7.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.886"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "2.196.63"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.997"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.951"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.32"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.420"
$ws.Range("E27").Value = "  +4.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.40"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("E29").Value = "  -11.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.575"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.348"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.487"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.273"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05330"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.314"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7642"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01997"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.845"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.51"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.614"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4586"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.040"
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8488"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.25"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.883"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.512"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.21"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4221"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.520"
$ws.Range("E51").Value = "  -0.32%  "
